$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; this shifts existing rows 8..29 down to 9..30
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44544
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100103
$ws.Cells.Item(8, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(8, 9).Value = 100103004
$ws.Cells.Item(8, 10).Value = "Durazno"
$ws.Cells.Item(8, 11).Value = "Polar King"
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 250
$ws.Cells.Item(8, 14).Value = 21000
$ws.Cells.Item(8, 15).Value = 22000
$ws.Cells.Item(8, 16).Value = 21500
$ws.Cells.Item(8, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(8, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(8, 19).Value = 1194
$ws.Cells.Item(8, 20).Value = 18
